# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Perejil" (Vega Modelo de Temuco) at row 430,
# shifting all subsequent rows down by one (dimension grows from A1:R537 to A1:R538).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 430; this pushes rows 430..537 down to 431..538.
$ws.Rows.Item(430).Insert()

# Populate the newly inserted row 430 with the new weekly record.
$ws.Range("A430").Value = 10
$ws.Range("B430").Value = "Vega Modelo de Temuco"
$ws.Range("C430").Value = "La Araucanía"
$ws.Range("D430").Value = 45211
$ws.Range("E430").Value = 9
$ws.Range("F430").Value = 100112044
$ws.Range("G430").Value = "Perejil"
$ws.Range("H430").Value = "Sin especificar"
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 50
$ws.Range("K430").Value = 5000
$ws.Range("L430").Value = 5000
$ws.Range("M430").Value = 5000
$ws.Range("N430").Value = "$/docena de atados (3 kilos)"
$ws.Range("O430").Value = "Provincia de Cautín"
$ws.Range("P430").Value = 1667
$ws.Range("Q430").Value = 3
$ws.Range("R430").Value = "Hortaliza"
